$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 3 - Bilan")
$ws.Range("G23").Value = 0.083333333333333301
